# "Drop in RMI script results for 3.0"
# - Remove the "Texas Notes" worksheet entirely.
# - Update the citation link on the About sheet to the new (working) URL.
# - Update the rebate-qualifying market-share-change data point from 4.95% to 7.4%.
# - Leave "About" as the active / selected sheet with the default A1 selection,
#   and restore the MSCdtRPbQL sheet's last-known selection (A2).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the "Texas Notes" sheet -------------------------------------
$wsNotes = $wb.Worksheets.Item("Texas Notes")
$wsNotes.Delete()

# --- About sheet: point the citation hyperlink at the new URL ----------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B6").Value = "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf"

# --- MSCdtRPbQL sheet: refresh the rebate market-share-change figure ---
$wsData = $wb.Worksheets.Item("MSCdtRPbQL")
$wsData.Range("C2").Value = 0.074
$wsData.Range("A2").Select()

# --- Restore "About" as the active sheet with a plain A1 selection -----
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
